$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# --- Fix cell kind (text placeholder <-> number) by copying format from stable reference cells ---
# C23/D23/F23/G23 = text "0" style; E23/H23 = text "***.*" style; I23/J23 = numeric style
$ws.Range("C23").Copy($ws.Range("C14"))
$ws.Range("I23").Copy($ws.Range("D14"))
$ws.Range("I23").Copy($ws.Range("E14"))
$ws.Range("I23").Copy($ws.Range("G14"))
$ws.Range("I23").Copy($ws.Range("H14"))
$ws.Range("C23").Copy($ws.Range("F15"))
$ws.Range("C23").Copy($ws.Range("C20"))
$ws.Range("C23").Copy($ws.Range("C22"))
$ws.Range("C23").Copy($ws.Range("F27"))
$ws.Range("I23").Copy($ws.Range("D29"))
$ws.Range("I23").Copy($ws.Range("E29"))
$ws.Range("I23").Copy($ws.Range("G29"))
$ws.Range("I23").Copy($ws.Range("H29"))
$ws.Range("I23").Copy($ws.Range("D30"))
$ws.Range("I23").Copy($ws.Range("E30"))
$ws.Range("I23").Copy($ws.Range("G30"))
$ws.Range("I23").Copy($ws.Range("H30"))

# --- Apply new values ---
$ws.Range("C14").Value = "0"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 100
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 40
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = "0"
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = 100
$ws.Range("M15").Value = 113.333333333333
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 240
$ws.Range("J16").Value = 217
$ws.Range("K16").Value = 10.599078341013
$ws.Range("L16").Value = 90.47619047619
$ws.Range("M16").Value = 81.818181818181
$ws.Range("N16").Value = -81.45285935085
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 93.75
$ws.Range("I17").Value = 314
$ws.Range("J17").Value = 272
$ws.Range("K17").Value = 15.441176470588
$ws.Range("L17").Value = 85.798816568047
$ws.Range("M17").Value = 155.284552845528
$ws.Range("N17").Value = -32.618025751073
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 6.25
$ws.Range("I18").Value = 183
$ws.Range("J18").Value = 182
$ws.Range("K18").Value = 0.54945054945
$ws.Range("L18").Value = 42.96875
$ws.Range("M18").Value = 3.389830508474
$ws.Range("N18").Value = -90.886454183266
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = -23.684210526315
$ws.Range("G19").Value = 169
$ws.Range("H19").Value = -21.89349112426
$ws.Range("I19").Value = 1640
$ws.Range("J19").Value = 1701
$ws.Range("K19").Value = -3.586125808348
$ws.Range("L19").Value = -5.52995391705
$ws.Range("M19").Value = 0.06101281269
$ws.Range("N19").Value = -76.668089344145
$ws.Range("C20").Value = "0"
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 100
$ws.Range("L20").Value = -32.051282051282
$ws.Range("N20").Value = -88.577586206896
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = -25.454545454545
$ws.Range("F21").Value = 198
$ws.Range("G21").Value = 228
$ws.Range("H21").Value = -13.157894736842
$ws.Range("I21").Value = 2469
$ws.Range("J21").Value = 2446
$ws.Range("K21").Value = 0.940310711365
$ws.Range("L21").Value = 9.782125389061
$ws.Range("M21").Value = 15.86109807602
$ws.Range("N21").Value = -78.140770252324
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 78
$ws.Range("K22").Value = -2.564102564102
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 26.666666666666
$ws.Range("C24").Value = 40
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").Value = 232
$ws.Range("H24").Value = 0.869565217391
$ws.Range("I24").Value = 2375
$ws.Range("J24").Value = 2908
$ws.Range("K24").Value = -18.328748280605
$ws.Range("L24").Value = -7.802795031055
$ws.Range("M24").Value = 29.498364231188
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = 45
$ws.Range("E25").Value = -35.555555555555
$ws.Range("G25").Value = 234
$ws.Range("H25").Value = -17.521367521367
$ws.Range("I25").Value = 2173
$ws.Range("J25").Value = 2762
$ws.Range("K25").Value = -21.325126719768
$ws.Range("L25").Value = -14.750882699097
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = 5.882352941176
$ws.Range("I26").Value = 690
$ws.Range("J26").Value = 697
$ws.Range("K26").Value = -1.004304160688
$ws.Range("L26").Value = 9.004739336492
$ws.Range("M26").Value = 59.353348729792
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = "0"
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = 41.666666666666
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 42.857142857142
$ws.Range("I28").Value = 119
$ws.Range("J28").Value = 110
$ws.Range("K28").Value = 8.181818181818
$ws.Range("L28").Value = 21.428571428571
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = 28.571428571428
$ws.Range("N29").Value = -30.76923076923
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -16.666666666666
$ws.Range("N30").Value = -50
$ws.Range("L31").Value = -33.333333333333
